$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update the workbook's reporting period
$ws.Name = "5-6-2018 to 10-6-2018"

# Row 8: hours changed from 2 to 4
$ws.Range("C8").Value = 4

# Row 10: new task entry (matches the green "Hands on" look used by B8/B9)
$ws.Range("A10").Value = 6
$ws.Range("B8").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").Value = "Reading  Large Scale Distributed Deep Networks"
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 43318

# Row 11: new task entry (keeps the default blue "Reading" look already on B11)
$ws.Range("A11").Value = 7
$ws.Range("B11").Value = "Tensorflow  Hands on mini batch stochastic gradeint descent (SGD)"
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 43318

# Update active cell selection
$ws.Range("I12").Select()
